$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Insert a new row at position 41 (shifts rows 41..48 down to 42..49)
$ws.Rows(41).Insert()

# New cell content for inserted row
$ws.Cells.Item(41, 1).Value = "juenger als 23 oder vor 1940 geboren?"

# Match formatting: A41 uses default style, B41 uses same style as B40 (integer number format)
$ws.Cells.Item(41, 2).Style = $ws.Cells.Item(40, 2).Style
$ws.Cells.Item(41, 2).NumberFormat = $ws.Cells.Item(40, 2).NumberFormat

# Update selection to match new state
$ws.Range("B39").Select() | Out-Null
